# post week 8 updates
#
# 1. New "Abstract Title" paragraph style (based on Normal, followed by
#    Abstract) - centered, bold, small, dark-blue heading that introduces
#    the Abstract block.
# 2. "Abstract" style now opens with less space before it (300 -> 100
#    twips) since the new Abstract Title supplies the space above.
# 3. New "Footnote Block Text" paragraph style (based on Footnote Text,
#    followed by Footnote Text) mirroring the existing "Block Text" style
#    but for footnotes - indented block quote look inside footnotes.

$d = $word.ActiveDocument

# --- 1. Abstract Title -----------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060   # RGB 34 5A 8A packed as BGR (0x8A5A34)

# --- 2. Abstract: tighten the space above the paragraph --------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. Footnote Block Text -------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "Added Abstract Title + Footnote Block Text styles; retuned Abstract spacing"
